$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.233.23'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '2.514.63'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''313.72'
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('D6').Value = '''93.25'
$ws.Range('E6').Value = '  -6.13%  '
$ws.Range('D7').Value = '''0.569'
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.526'
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('D10').Value = '''35.40'
$ws.Range('E10').Value = '  -5.10%  '
$ws.Range('D11').Value = '''0.0800'
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('D12').Value = '''7.54'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '2.901.79'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').Value = '2.510.96'
$ws.Range('E15').Value = '  -2.58%  '
$ws.Range('D16').Value = '''15.25'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '''0.843'
$ws.Range('E17').Value = '  -4.01%  '
$ws.Range('D18').Value = '42.369.69'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').Value = '''12.74'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('D20').Value = '''6.50'
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('D21').Value = '0.0₃0951'
$ws.Range('E21').Value = '  -4.12%  '
$ws.Range('D22').Value = '''70.32'
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').Value = '''248.62'
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('D24').Value = '''2.92'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').Value = '''1.99'
$ws.Range('E25').Value = '  -3.75%  '
$ws.Range('D26').Value = '''26.19'
$ws.Range('E26').Value = '  -5.66%  '
$ws.Range('D27').Value = '''0.997'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').Value = '''2.36'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.06'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''38.51'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '''5.84'
$ws.Range('E31').Value = '  -5.91%  '
$ws.Range('D32').Value = '''155.51'
$ws.Range('E32').Value = '  -1.60%  '
$ws.Range('D33').Value = '''19.30'
$ws.Range('E33').Value = '  +2.36%  '
$ws.Range('D34').Value = '''2.09'
$ws.Range('E34').Value = '  -2.23%  '
$ws.Range('D35').Value = '''3.25'
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''2.62'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.0775'
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('D38').Value = '''0.109'
$ws.Range('E38').Value = '  -5.64%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '''0.118'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '''23.72'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').Value = '''2.33'
$ws.Range('E41').Value = '  +11.80%  '
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').Value = '''3.75'
$ws.Range('E43').Value = '  -4.17%  '
$ws.Range('D44').Value = '''3.27'
$ws.Range('E44').Value = '  -5.68%  '
$ws.Range('D45').Value = '''0.0297'
$ws.Range('E45').Value = '  -2.79%  '
$ws.Range('D46').Value = '2.009.20'
$ws.Range('E46').Value = '  -3.00%  '
$ws.Range('D47').Value = '''83.87'
$ws.Range('E47').Value = '  -2.93%  '
$ws.Range('D48').Value = '''8.76'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('D49').Value = '2.759.21'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('D50').Value = '''72.21'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('D51').Value = '''101.34'
$ws.Range('E51').Value = '  -2.16%  '
